# Scheduled market-data refresh
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ (columns H-N) for the
# leve rows whose underlying Universalis market data changed since the last run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 369.95  # ALC!H107: 355.42856 -> 369.95
$ws.Cells.Item(107, 9).Value = 369.95  # ALC!I107: 355.42856 -> 369.95
$ws.Cells.Item(107, 11).Value = 369.95  # ALC!K107: 355.42856 -> 369.95
$ws.Cells.Item(107, 13).Value = 1550.05  # ALC!M107: 1564.57144 -> 1550.05

$ws.Cells.Item(113, 8).Value = 45458388  # ALC!H113: 50003930 -> 45458388
$ws.Cells.Item(113, 10).Value = 66670570  # ALC!J113: 80004080 -> 66670570
$ws.Cells.Item(113, 12).Value = 66670570  # ALC!L113: 80004080 -> 66670570
$ws.Cells.Item(113, 14).Value = -66677078  # ALC!N113: -80010588 -> -66677078

$ws.Cells.Item(118, 8).Value = 315.7  # ALC!H118: 330 -> 315.7
$ws.Cells.Item(118, 9).Value = 296  # ALC!I118: 309.625 -> 296
$ws.Cells.Item(118, 11).Value = 888  # ALC!K118: 928.875 -> 888
$ws.Cells.Item(118, 13).Value = 769  # ALC!M118: 728.125 -> 769

$ws.Cells.Item(127, 8).Value = 4079.2  # ALC!H127: 2229.3 -> 4079.2
$ws.Cells.Item(127, 9).Value = 698.6667  # ALC!I127: 427.57144 -> 698.6667
$ws.Cells.Item(127, 10).Value = 9150  # ALC!J127: 6433.3335 -> 9150
$ws.Cells.Item(127, 11).Value = 2096.0001  # ALC!K127: 1282.71432 -> 2096.0001
$ws.Cells.Item(127, 12).Value = 27450  # ALC!L127: 19300.0005 -> 27450
$ws.Cells.Item(127, 13).Value = 2863.9999  # ALC!M127: 3677.28568 -> 2863.9999
$ws.Cells.Item(127, 14).Value = -37370  # ALC!N127: -29220.0005 -> -37370

$ws.Cells.Item(132, 8).Value = 2177.6924  # ALC!H132: 2400.25 -> 2177.6924
$ws.Cells.Item(132, 9).Value = 1848.6364  # ALC!I132: 1984.909 -> 1848.6364
$ws.Cells.Item(132, 10).Value = 3987.5  # ALC!J132: 6969 -> 3987.5
$ws.Cells.Item(132, 11).Value = 5545.9092  # ALC!K132: 5954.727000000001 -> 5545.9092
$ws.Cells.Item(132, 12).Value = 11962.5  # ALC!L132: 20907 -> 11962.5
$ws.Cells.Item(132, 13).Value = -3015.9092  # ALC!M132: -3424.727000000001 -> -3015.9092
$ws.Cells.Item(132, 14).Value = -17022.5  # ALC!N132: -25967 -> -17022.5

$ws.Cells.Item(135, 8).Value = 2863.9412  # ALC!H135: 2759.611 -> 2863.9412
$ws.Cells.Item(135, 9).Value = 2692  # ALC!I135: 2578.2666 -> 2692
$ws.Cells.Item(135, 11).Value = 24228  # ALC!K135: 23204.3994 -> 24228
$ws.Cells.Item(135, 13).Value = -21693  # ALC!M135: -20669.3994 -> -21693

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 21675364  # ARM!H74: 20008174 -> 21675364
$ws.Cells.Item(74, 10).Value = 1259296.9  # ARM!J74: 1119586.1 -> 1259296.9
$ws.Cells.Item(74, 12).Value = 1259296.9  # ARM!L74: 1119586.1 -> 1259296.9
$ws.Cells.Item(74, 14).Value = -1261044.9  # ARM!N74: -1121334.1 -> -1261044.9

$ws.Cells.Item(77, 8).Value = 21675364  # ARM!H77: 20008174 -> 21675364
$ws.Cells.Item(77, 10).Value = 1259296.9  # ARM!J77: 1119586.1 -> 1259296.9
$ws.Cells.Item(77, 12).Value = 6296484.5  # ARM!L77: 5597930.5 -> 6296484.5
$ws.Cells.Item(77, 14).Value = -6305220.5  # ARM!N77: -5606666.5 -> -6305220.5

$ws.Cells.Item(92, 8).Value = 77446.75  # ARM!H92: 69996.336 -> 77446.75
$ws.Cells.Item(92, 10).Value = 79929  # ARM!J92: 69994.5 -> 79929
$ws.Cells.Item(92, 12).Value = 79929  # ARM!L92: 69994.5 -> 79929
$ws.Cells.Item(92, 14).Value = -84921  # ARM!N92: -74986.5 -> -84921

$ws.Cells.Item(110, 8).Value = 2927.2942  # ARM!H110: 3017.25 -> 2927.2942
$ws.Cells.Item(110, 9).Value = 2494.6155  # ARM!I110: 2578.5 -> 2494.6155
$ws.Cells.Item(110, 11).Value = 2494.6155  # ARM!K110: 2578.5 -> 2494.6155
$ws.Cells.Item(110, 13).Value = -449.6154999999999  # ARM!M110: -533.5 -> -449.6154999999999

$ws.Cells.Item(132, 8).Value = 6323.7427  # ARM!H132: 6214.75 -> 6323.7427
$ws.Cells.Item(132, 9).Value = 4216.64  # ARM!I132: 4146.769 -> 4216.64
$ws.Cells.Item(132, 11).Value = 12649.92  # ARM!K132: 12440.307 -> 12649.92
$ws.Cells.Item(132, 13).Value = -10119.92  # ARM!M132: -9910.307000000001 -> -10119.92

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 2707.7144  # BSM!H107: 3664.5 -> 2707.7144
$ws.Cells.Item(107, 9).Value = 3181.2  # BSM!I107: 5805 -> 3181.2
$ws.Cells.Item(107, 11).Value = 3181.2  # BSM!K107: 5805 -> 3181.2
$ws.Cells.Item(107, 13).Value = -1261.2  # BSM!M107: -3885 -> -1261.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 901670.8  # CRP!H31: 710973.0600000001 -> 901670.8
$ws.Cells.Item(31, 9).Value = 2829.3635  # CRP!I31: 2705.5 -> 2829.3635
$ws.Cells.Item(31, 10).Value = 1560821.2  # CRP!J31: 1377577.9 -> 1560821.2
$ws.Cells.Item(31, 11).Value = 2829.3635  # CRP!K31: 2705.5 -> 2829.3635
$ws.Cells.Item(31, 12).Value = 1560821.2  # CRP!L31: 1377577.9 -> 1560821.2
$ws.Cells.Item(31, 13).Value = -2534.3635  # CRP!M31: -2410.5 -> -2534.3635
$ws.Cells.Item(31, 14).Value = -1561411.2  # CRP!N31: -1378167.9 -> -1561411.2

$ws.Cells.Item(34, 8).Value = 901670.8  # CRP!H34: 710973.0600000001 -> 901670.8
$ws.Cells.Item(34, 9).Value = 2829.3635  # CRP!I34: 2705.5 -> 2829.3635
$ws.Cells.Item(34, 10).Value = 1560821.2  # CRP!J34: 1377577.9 -> 1560821.2
$ws.Cells.Item(34, 11).Value = 2829.3635  # CRP!K34: 2705.5 -> 2829.3635
$ws.Cells.Item(34, 12).Value = 1560821.2  # CRP!L34: 1377577.9 -> 1560821.2
$ws.Cells.Item(34, 13).Value = -2627.3635  # CRP!M34: -2503.5 -> -2627.3635
$ws.Cells.Item(34, 14).Value = -1561225.2  # CRP!N34: -1377981.9 -> -1561225.2

$ws.Cells.Item(70, 8).Value = 40000  # CRP!H70: 24900 -> 40000
$ws.Cells.Item(70, 10).Value = 40000  # CRP!J70: 24900 -> 40000
$ws.Cells.Item(70, 12).Value = 40000  # CRP!L70: 24900 -> 40000
$ws.Cells.Item(70, 14).Value = -40630  # CRP!N70: -25530 -> -40630

$ws.Cells.Item(73, 8).Value = 40000  # CRP!H73: 24900 -> 40000
$ws.Cells.Item(73, 10).Value = 40000  # CRP!J73: 24900 -> 40000
$ws.Cells.Item(73, 12).Value = 40000  # CRP!L73: 24900 -> 40000
$ws.Cells.Item(73, 14).Value = -42184  # CRP!N73: -27084 -> -42184

$ws.Cells.Item(99, 8).Value = 4012  # CRP!H99: 3205.2 -> 4012
$ws.Cells.Item(99, 9).Value = 4012  # CRP!I99: 3205.2 -> 4012
$ws.Cells.Item(99, 11).Value = 4012  # CRP!K99: 3205.2 -> 4012
$ws.Cells.Item(99, 13).Value = -2514  # CRP!M99: -1707.2 -> -2514

$ws.Cells.Item(107, 8).Value = 2426.3333  # CRP!H107: 2346 -> 2426.3333
$ws.Cells.Item(107, 9).Value = 957.3333  # CRP!I107: 943 -> 957.3333
$ws.Cells.Item(107, 11).Value = 957.3333  # CRP!K107: 943 -> 957.3333
$ws.Cells.Item(107, 13).Value = 962.6667  # CRP!M107: 977 -> 962.6667

$ws.Cells.Item(126, 8).Value = 4012  # CRP!H126: 3205.2 -> 4012
$ws.Cells.Item(126, 9).Value = 4012  # CRP!I126: 3205.2 -> 4012
$ws.Cells.Item(126, 11).Value = 12036  # CRP!K126: 9615.599999999999 -> 12036
$ws.Cells.Item(126, 13).Value = -9566  # CRP!M126: -7145.599999999999 -> -9566

$ws.Cells.Item(134, 8).Value = 6892.4  # CRP!H134: 4501 -> 6892.4
$ws.Cells.Item(134, 9).Value = 3678.3333  # CRP!I134: 1862.8 -> 3678.3333
$ws.Cells.Item(134, 11).Value = 11034.9999  # CRP!K134: 5588.4 -> 11034.9999
$ws.Cells.Item(134, 13).Value = -8499.999899999999  # CRP!M134: -3053.4 -> -8499.999899999999

$ws.Cells.Item(141, 8).Value = 213757.67  # CRP!H141: 229372.1 -> 213757.67
$ws.Cells.Item(141, 10).Value = 224190.19  # CRP!J141: 242409.3 -> 224190.19
$ws.Cells.Item(141, 12).Value = 224190.19  # CRP!L141: 242409.3 -> 224190.19
$ws.Cells.Item(141, 14).Value = -234550.19  # CRP!N141: -252769.3 -> -234550.19

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 20835716  # CUL!H4: 20835760 -> 20835716
$ws.Cells.Item(4, 9).Value = 12100207  # CUL!I4: 12100267 -> 12100207
$ws.Cells.Item(4, 11).Value = 36300621  # CUL!K4: 36300801 -> 36300621
$ws.Cells.Item(4, 13).Value = -36300509  # CUL!M4: -36300689 -> -36300509

$ws.Cells.Item(5, 8).Value = 1915.1875  # CUL!H5: 1740.6 -> 1915.1875
$ws.Cells.Item(5, 9).Value = 1967.4286  # CUL!I5: 1867.6666 -> 1967.4286
$ws.Cells.Item(5, 10).Value = 1549.5  # CUL!J5: 1359.4 -> 1549.5
$ws.Cells.Item(5, 11).Value = 5902.2858  # CUL!K5: 5602.9998 -> 5902.2858
$ws.Cells.Item(5, 12).Value = 4648.5  # CUL!L5: 4078.2 -> 4648.5
$ws.Cells.Item(5, 13).Value = -5790.2858  # CUL!M5: -5490.9998 -> -5790.2858
$ws.Cells.Item(5, 14).Value = -4872.5  # CUL!N5: -4302.200000000001 -> -4872.5

$ws.Cells.Item(6, 8).Value = 915.8333  # CUL!H6: 958.8 -> 915.8333
$ws.Cells.Item(6, 9).Value = 899  # CUL!I6: 948.5 -> 899
$ws.Cells.Item(6, 11).Value = 2697  # CUL!K6: 2845.5 -> 2697
$ws.Cells.Item(6, 13).Value = -2584  # CUL!M6: -2732.5 -> -2584

$ws.Cells.Item(7, 8).Value = 1031.5  # CUL!H7: 1668 -> 1031.5
$ws.Cells.Item(7, 9).Value = 500.66666  # CUL!I7: 501 -> 500.66666
$ws.Cells.Item(7, 10).Value = 1562.3334  # CUL!J7: 4002 -> 1562.3334
$ws.Cells.Item(7, 11).Value = 1501.99998  # CUL!K7: 1503 -> 1501.99998
$ws.Cells.Item(7, 12).Value = 4687.0002  # CUL!L7: 12006 -> 4687.0002
$ws.Cells.Item(7, 13).Value = -1389.99998  # CUL!M7: -1391 -> -1389.99998
$ws.Cells.Item(7, 14).Value = -4911.0002  # CUL!N7: -12230 -> -4911.0002

$ws.Cells.Item(10, 8).Value = 155.25  # CUL!H10: 142.8 -> 155.25
$ws.Cells.Item(10, 9).Value = 39.666668  # CUL!I10: 53 -> 39.666668
$ws.Cells.Item(10, 11).Value = 119.000004  # CUL!K10: 159 -> 119.000004
$ws.Cells.Item(10, 13).Value = 19.999996  # CUL!M10: -20 -> 19.999996

$ws.Cells.Item(11, 8).Value = 141.77777  # CUL!H11: 130.6 -> 141.77777
$ws.Cells.Item(11, 9).Value = 141.77777  # CUL!I11: 130.6 -> 141.77777
$ws.Cells.Item(11, 11).Value = 425.33331  # CUL!K11: 391.8 -> 425.33331
$ws.Cells.Item(11, 13).Value = -285.33331  # CUL!M11: -251.8 -> -285.33331

$ws.Cells.Item(13, 8).Value = 5356.6665  # CUL!H13: 5761.4 -> 5356.6665
$ws.Cells.Item(13, 10).Value = 5556  # CUL!J13: 6667.5 -> 5556
$ws.Cells.Item(13, 12).Value = 16668  # CUL!L13: 20002.5 -> 16668
$ws.Cells.Item(13, 14).Value = -17004  # CUL!N13: -20338.5 -> -17004

$ws.Cells.Item(50, 8).Value = 467.21738  # CUL!H50: 472.54544 -> 467.21738
$ws.Cells.Item(50, 9).Value = 349.5  # CUL!I50: 349 -> 349.5
$ws.Cells.Item(50, 11).Value = 1048.5  # CUL!K50: 1047 -> 1048.5
$ws.Cells.Item(50, 13).Value = -567.5  # CUL!M50: -566 -> -567.5

$ws.Cells.Item(53, 8).Value = 467.21738  # CUL!H53: 472.54544 -> 467.21738
$ws.Cells.Item(53, 9).Value = 349.5  # CUL!I53: 349 -> 349.5
$ws.Cells.Item(53, 11).Value = 1048.5  # CUL!K53: 1047 -> 1048.5
$ws.Cells.Item(53, 13).Value = -567.5  # CUL!M53: -566 -> -567.5

$ws.Cells.Item(122, 8).Value = 1185.35  # CUL!H122: 1210.8948 -> 1185.35
$ws.Cells.Item(122, 10).Value = 1609.0714  # CUL!J122: 1679 -> 1609.0714
$ws.Cells.Item(122, 12).Value = 14481.6426  # CUL!L122: 15111 -> 14481.6426
$ws.Cells.Item(122, 14).Value = -19381.6426  # CUL!N122: -20011 -> -19381.6426

$ws.Cells.Item(129, 8).Value = 15875813  # CUL!H129: 18521702 -> 15875813
$ws.Cells.Item(129, 9).Value = 3455  # CUL!I129: 4105.8 -> 3455
$ws.Cells.Item(129, 10).Value = 22224756  # CUL!J129: 25643854 -> 22224756
$ws.Cells.Item(129, 11).Value = 10365  # CUL!K129: 12317.4 -> 10365
$ws.Cells.Item(129, 12).Value = 66674268  # CUL!L129: 76931562 -> 66674268
$ws.Cells.Item(129, 13).Value = -5365  # CUL!M129: -7317.400000000001 -> -5365
$ws.Cells.Item(129, 14).Value = -66684268  # CUL!N129: -76941562 -> -66684268

$ws.Cells.Item(131, 8).Value = 6821.1284  # CUL!H131: 5694.163 -> 6821.1284
$ws.Cells.Item(131, 10).Value = 7273.276  # CUL!J131: 5741.41 -> 7273.276
$ws.Cells.Item(131, 12).Value = 21819.828  # CUL!L131: 17224.23 -> 21819.828
$ws.Cells.Item(131, 14).Value = -31899.828  # CUL!N131: -27304.23 -> -31899.828

$ws.Cells.Item(135, 8).Value = 1915.1875  # CUL!H135: 1740.6 -> 1915.1875
$ws.Cells.Item(135, 9).Value = 1967.4286  # CUL!I135: 1867.6666 -> 1967.4286
$ws.Cells.Item(135, 10).Value = 1549.5  # CUL!J135: 1359.4 -> 1549.5
$ws.Cells.Item(135, 11).Value = 17706.8574  # CUL!K135: 16808.9994 -> 17706.8574
$ws.Cells.Item(135, 12).Value = 13945.5  # CUL!L135: 12234.6 -> 13945.5
$ws.Cells.Item(135, 13).Value = -15171.8574  # CUL!M135: -14273.9994 -> -15171.8574
$ws.Cells.Item(135, 14).Value = -19015.5  # CUL!N135: -17304.6 -> -19015.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 392.8889  # GSM!H107: 410 -> 392.8889
$ws.Cells.Item(107, 9).Value = 357.06668  # GSM!I107: 375.2857 -> 357.06668
$ws.Cells.Item(107, 11).Value = 357.06668  # GSM!K107: 375.2857 -> 357.06668
$ws.Cells.Item(107, 13).Value = 1562.93332  # GSM!M107: 1544.7143 -> 1562.93332

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(123, 8).Value = 79326.336  # LTW!H123: 79326.664 -> 79326.336
$ws.Cells.Item(123, 10).Value = 79326.336  # LTW!J123: 79326.664 -> 79326.336
$ws.Cells.Item(123, 12).Value = 79326.336  # LTW!L123: 79326.664 -> 79326.336
$ws.Cells.Item(123, 14).Value = -89126.336  # LTW!N123: -89126.664 -> -89126.336

$ws.Cells.Item(132, 8).Value = 681614.9  # LTW!H132: 568412.3 -> 681614.9
$ws.Cells.Item(132, 9).Value = 17601.5  # LTW!I132: 14561.134 -> 17601.5
$ws.Cells.Item(132, 11).Value = 52804.5  # LTW!K132: 43683.402 -> 52804.5
$ws.Cells.Item(132, 13).Value = -50274.5  # LTW!M132: -41153.402 -> -50274.5

$ws.Cells.Item(135, 8).Value = 0  # LTW!H135: 50000 -> 0
$ws.Cells.Item(135, 10).Value = 0  # LTW!J135: 50000 -> 0
$ws.Cells.Item(135, 12).Value = 0  # LTW!L135: 50000 -> 0
$ws.Cells.Item(135, 14).ClearContents()  # LTW!N135: removed (was -60140)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 10975  # WVR!H41: 13000 -> 10975
$ws.Cells.Item(41, 10).Value = 10950  # WVR!J41: 15000 -> 10950
$ws.Cells.Item(41, 12).Value = 10950  # WVR!L41: 15000 -> 10950
$ws.Cells.Item(41, 14).Value = -11730  # WVR!N41: -15780 -> -11730

$ws.Cells.Item(126, 8).Value = 4500.125  # WVR!H126: 3714.5715 -> 4500.125
$ws.Cells.Item(126, 9).Value = 4800.4  # WVR!I126: 3500.75 -> 4800.4
$ws.Cells.Item(126, 11).Value = 14401.2  # WVR!K126: 10502.25 -> 14401.2
$ws.Cells.Item(126, 13).Value = -11931.2  # WVR!M126: -8032.25 -> -11931.2

$ws.Cells.Item(136, 8).Value = 4679.4116  # WVR!H136: 4270.316 -> 4679.4116
$ws.Cells.Item(136, 9).Value = 4596.875  # WVR!I136: 4174.222 -> 4596.875
$ws.Cells.Item(136, 11).Value = 13790.625  # WVR!K136: 12522.666 -> 13790.625
$ws.Cells.Item(136, 13).Value = -11240.625  # WVR!M136: -9972.665999999999 -> -11240.625
